# Elimina participante desde panel de administración
# Removes the row for "Harrison Driver_20251202_130401" (the "Sin SmartScore"
# duplicate submission made at 13:04:01), shifting all subsequent rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetId = "Harrison Driver_20251202_130401"
$lastRow = $ws.UsedRange.Rows.Count
$rowToDelete = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $val = $ws.Cells.Item($r, 1).Value2
    if ($val -eq $targetId) {
        $rowToDelete = $r
        break
    }
}

if ($rowToDelete -gt 0) {
    $ws.Rows.Item($rowToDelete).Delete()
}
